$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H6").Formula = "=C6-C5"
$ws.Range("G7").Formula = "=C7-C5"
$ws.Range("H12").Formula = "=C12-C11"
$ws.Range("G13").Formula = "=C13-C11"
$ws.Range("H18").Formula = "=C18-C17"
$ws.Range("G19").Formula = "=C19-C17"
$ws.Range("H24").Formula = "=C24-C23"
$ws.Range("G25").Formula = "=C25-C23"
$ws.Range("H30").Formula = "=C30-C29"
$ws.Range("G31").Formula = "=C31-C29"
$ws.Range("H36").Formula = "=C36-C35"
$ws.Range("G37").Formula = "=C37-C35"
$ws.Range("G39").Formula = "=AVERAGE(G7:G37)"
$ws.Range("H39").Formula = "=AVERAGE(H7:H37)"
$ws.Range("H39").NumberFormat = "0.00"

[void]$ws.Range("H36").Select()
